# Update Palorsennan Train Lines.xlsx
# - Adds two new train lines (The Capital Line, The Lake Line)
# - Updates Purpose/Tracks values for a few existing lines
# - Fixes a stray space in "Snubavik- Bofoker" -> "Snubavik-Bofoker"
# - Re-sorts the lines table alphabetically by Line Name
# - Shifts the Companies/Ownerships/Description table down to make room

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start clean: wipe all existing cell content and per-cell formatting
# (column widths/default column styles are unaffected).
$ws.UsedRange.Clear()

# ---- Header row ----
$ws.Cells.Item(1, 1).Value = "Line Name"
$ws.Cells.Item(1, 2).Value = "Ownership"
$ws.Cells.Item(1, 3).Value = "Tracks"
$ws.Cells.Item(1, 4).Value = "Purpose"
$ws.Cells.Item(1, 5).Value = "Line Ends"
$ws.Cells.Item(1, 6).Value = "Stops"
$ws.Cells.Item(1, 7).Value = "Branches"
$ws.Range("A1:G1").Font.Bold = $true

# ---- Train line rows (2-10), alphabetical by Line Name ----

# Row 2: The Capital Line
$ws.Cells.Item(2, 1).Value = "The Capital Line"
$ws.Cells.Item(2, 2).Value = "Palorsenna Rail"
$ws.Cells.Item(2, 3).Value = "4 ("
$ws.Cells.Item(2, 5).Value = "Raunahild-Raupige"
$ws.Cells.Item(2, 6).Value = "Raunahild-Rivervo-Raupige"

# Row 3: The Coastal Line
$ws.Cells.Item(3, 1).Value = "The Coastal Line"
$ws.Cells.Item(3, 2).Value = "Palorsenna Rail"
$ws.Cells.Item(3, 3).Value = "4 ("
$ws.Cells.Item(3, 5).Value = "Raunahild-Oyagerdur"

# Row 4: The Eastern Line
$ws.Cells.Item(4, 1).Value = "The Eastern Line"
$ws.Cells.Item(4, 2).Value = "Palorsenna Rail"
$ws.Cells.Item(4, 3).Value = "4 ("
$ws.Cells.Item(4, 5).Value = "Donland-Solder"

# Row 5: The Industrial Line
$ws.Cells.Item(5, 1).Value = "The Industrial Line"
$ws.Cells.Item(5, 2).Value = "Palorsenna Rail"
$ws.Cells.Item(5, 3).Value = "4 ("
$ws.Cells.Item(5, 5).Value = "Snubavik-Bofoker"
$ws.Cells.Item(5, 7).Value = "Bofoker-Grimsjahver"

# Row 6: The Lake Line
$ws.Cells.Item(6, 1).Value = "The Lake Line"
$ws.Cells.Item(6, 2).Value = "Palorsenna Rail"
$ws.Cells.Item(6, 3).Value = "4 ("
$ws.Cells.Item(6, 5).Value = "Raupige-Raupige"
$ws.Cells.Item(6, 6).Value = "Raupige-Snubavik-Karldalla-Raupige"

# Row 7: The Prosperity Line
$ws.Cells.Item(7, 1).Value = "The Prosperity Line"
$ws.Cells.Item(7, 2).Value = "The Prosperity Line Company"
$ws.Cells.Item(7, 3).Value = "4 (2 Maglev, 2 Low Speed)"
$ws.Cells.Item(7, 4).Value = "Very High Speed Passenger Transportation, Freight"
$ws.Cells.Item(7, 5).Value = "Raunahild-Montara"
$ws.Cells.Item(7, 6).Value = "Raunahild-Raupige-Snubavik-Hvolstad-Montara"

# Row 8: The Riverlands Line
$ws.Cells.Item(8, 1).Value = "The Riverlands Line"
$ws.Cells.Item(8, 2).Value = "Palorsenna Rail"
$ws.Cells.Item(8, 3).Value = "2 High Speed"
$ws.Cells.Item(8, 4).Value = "High Speed Passenger Transportation"
$ws.Cells.Item(8, 5).Value = "Snubavik-Bofoker"
$ws.Cells.Item(8, 6).Value = "Snubavik-Gravden-Lokernes-Kjersnes-Barbakvik-Laufargar-Bofoker"

# Row 9: The Sulion Line
$ws.Cells.Item(9, 1).Value = "The Sulion Line"
$ws.Cells.Item(9, 2).Value = "Palorsenna Rail"
$ws.Cells.Item(9, 3).Value = "2 Low Speed"
$ws.Cells.Item(9, 4).Value = "Low Speed Passenger Transportation, Freight"
$ws.Cells.Item(9, 5).Value = "Valvegils-Onesos"
$ws.Cells.Item(9, 6).Value = "Valvegils-Sokervipo-Snubageid-Kateydanes-Onesos"
$ws.Cells.Item(9, 7).Value = "Snubageid-Redcke"

# Row 10: The Western Line
$ws.Cells.Item(10, 1).Value = "The Western Line"
$ws.Cells.Item(10, 2).Value = "Palorsenna Rail"
$ws.Cells.Item(10, 3).Value = "4 (2 Medium Speed, 2 Low Speed)"
$ws.Cells.Item(10, 4).Value = "Medium Speed Passenger Transportation, Frieght"
$ws.Cells.Item(10, 5).Value = "Raunahild-Mazion"
$ws.Cells.Item(10, 6).Value = "Raunahild-Horestad-Stolafsf-Raudagast-Onesos-Mazion"
$ws.Cells.Item(10, 7).Value = "Onesos-Arsycharann-Mazion"

# ---- Companies / Ownerships / Description table (rows 13-15) ----
$ws.Cells.Item(13, 2).Value = "Companies"
$ws.Cells.Item(13, 3).Value = "Ownerships"
$ws.Cells.Item(13, 4).Value = "Description"
$ws.Range("B13:D13").Font.Bold = $true

$ws.Cells.Item(14, 2).Value = "Palorsenna Rail"
$ws.Cells.Item(14, 3).Value = "Palorsennan Government"
$ws.Cells.Item(14, 4).Value = "Public utility rail."

$ws.Cells.Item(15, 2).Value = "The Prosperity Line Company"
$ws.Cells.Item(15, 3).Value = "Palorsenna Rail, (Artizore Rail)"
$ws.Cells.Item(15, 4).Value = "Joint company for the running of the prosperity line."

# ---- Column F got wider to fit the new, longer stop lists ----
$ws.Columns.Item(6).ColumnWidth = 67.8

# ---- Selection ends on B10 ----
$ws.Range("B10").Select() | Out-Null
